$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 98
$ws_ALC.Range("H98").Value = 1181.0667
$ws_ALC.Range("I98").Value = 1055.0769
$ws_ALC.Range("J98").Value = 2000
$ws_ALC.Range("K98").Value = 1055.0769
$ws_ALC.Range("L98").Value = 2000
$ws_ALC.Range("M98").Value = 442.9231
$ws_ALC.Range("N98").Value = -4996

# ALC row 122
$ws_ALC.Range("H122").Value = 1181.0667
$ws_ALC.Range("I122").Value = 1055.0769
$ws_ALC.Range("J122").Value = 2000
$ws_ALC.Range("K122").Value = 3165.2307
$ws_ALC.Range("L122").Value = 6000
$ws_ALC.Range("M122").Value = -715.2307000000001
$ws_ALC.Range("N122").Value = -10900

# ALC row 135
$ws_ALC.Range("H135").Value = 321.39215
$ws_ALC.Range("I135").Value = 288.42554
$ws_ALC.Range("J135").Value = 708.75
$ws_ALC.Range("K135").Value = 2595.82986
$ws_ALC.Range("L135").Value = 6378.75
$ws_ALC.Range("M135").Value = -60.82986000000028
$ws_ALC.Range("N135").Value = -11448.75

# ALC row 138
$ws_ALC.Range("H138").Value = 1493.3969
$ws_ALC.Range("I138").Value = 859.88
$ws_ALC.Range("J138").Value = 3930
$ws_ALC.Range("K138").Value = 2579.64
$ws_ALC.Range("L138").Value = 11790
$ws_ALC.Range("M138").Value = 2560.36
$ws_ALC.Range("N138").Value = -22070

# ARM row 74
$ws_ARM.Range("H74").Value = 4323.7
$ws_ARM.Range("I74").Value = 986.25
$ws_ARM.Range("J74").Value = 17673.5
$ws_ARM.Range("K74").Value = 986.25
$ws_ARM.Range("L74").Value = 17673.5
$ws_ARM.Range("M74").Value = -112.25
$ws_ARM.Range("N74").Value = -19421.5

# ARM row 77
$ws_ARM.Range("H77").Value = 4323.7
$ws_ARM.Range("I77").Value = 986.25
$ws_ARM.Range("J77").Value = 17673.5
$ws_ARM.Range("K77").Value = 4931.25
$ws_ARM.Range("L77").Value = 88367.5
$ws_ARM.Range("M77").Value = -563.25
$ws_ARM.Range("N77").Value = -97103.5

# ARM row 132
$ws_ARM.Range("H132").Value = 12575.904
$ws_ARM.Range("I132").Value = 12414
$ws_ARM.Range("J132").Value = 12723.091
$ws_ARM.Range("K132").Value = 37242
$ws_ARM.Range("L132").Value = 38169.273
$ws_ARM.Range("M132").Value = -34712
$ws_ARM.Range("N132").Value = -43229.273

# BSM row 52
$ws_BSM.Range("H52").Value = 0
$ws_BSM.Range("I52").Value = 0
$ws_BSM.Range("J52").Value = 0
$ws_BSM.Range("K52").Value = 0
$ws_BSM.Range("L52").Value = 0
$ws_BSM.Range("N52").ClearContents()

# BSM row 82
$ws_BSM.Range("H82").Value = 13480.429
$ws_BSM.Range("I82").Value = 10727.167
$ws_BSM.Range("J82").Value = 30000
$ws_BSM.Range("K82").Value = 10727.167
$ws_BSM.Range("L82").Value = 30000
$ws_BSM.Range("M82").Value = -10344.167
$ws_BSM.Range("N82").Value = -30766

# BSM row 85
$ws_BSM.Range("H85").Value = 13480.429
$ws_BSM.Range("I85").Value = 10727.167
$ws_BSM.Range("J85").Value = 30000
$ws_BSM.Range("K85").Value = 10727.167
$ws_BSM.Range("L85").Value = 30000
$ws_BSM.Range("M85").Value = -9401.166999999999
$ws_BSM.Range("N85").Value = -32652

# BSM row 116
$ws_BSM.Range("H116").Value = 0
$ws_BSM.Range("I116").Value = 0
$ws_BSM.Range("J116").Value = 0
$ws_BSM.Range("K116").Value = 0
$ws_BSM.Range("L116").Value = 0
$ws_BSM.Range("N116").ClearContents()

# BSM row 117
$ws_BSM.Range("H117").Value = 0
$ws_BSM.Range("I117").Value = 0
$ws_BSM.Range("J117").Value = 0
$ws_BSM.Range("K117").Value = 0
$ws_BSM.Range("L117").Value = 0
$ws_BSM.Range("N117").ClearContents()

# BSM row 118
$ws_BSM.Range("H118").Value = 0
$ws_BSM.Range("I118").Value = 0
$ws_BSM.Range("J118").Value = 0
$ws_BSM.Range("K118").Value = 0
$ws_BSM.Range("L118").Value = 0
$ws_BSM.Range("N118").ClearContents()

# BSM row 119
$ws_BSM.Range("H119").Value = 0
$ws_BSM.Range("I119").Value = 0
$ws_BSM.Range("J119").Value = 0
$ws_BSM.Range("K119").Value = 0
$ws_BSM.Range("L119").Value = 0
$ws_BSM.Range("N119").ClearContents()

# BSM row 120
$ws_BSM.Range("H120").Value = 48000
$ws_BSM.Range("I120").Value = 0
$ws_BSM.Range("J120").Value = 48000
$ws_BSM.Range("K120").Value = 0
$ws_BSM.Range("L120").Value = 48000
$ws_BSM.Range("N120").Value = -57676

# BSM row 121
$ws_BSM.Range("H121").Value = 0
$ws_BSM.Range("I121").Value = 0
$ws_BSM.Range("J121").Value = 0
$ws_BSM.Range("K121").Value = 0
$ws_BSM.Range("L121").Value = 0
$ws_BSM.Range("N121").ClearContents()

# BSM row 130
$ws_BSM.Range("H130").Value = 47780
$ws_BSM.Range("I130").Value = 0
$ws_BSM.Range("J130").Value = 47780
$ws_BSM.Range("K130").Value = 0
$ws_BSM.Range("L130").Value = 47780
$ws_BSM.Range("N130").Value = -57820

# BSM row 134
$ws_BSM.Range("H134").Value = 1079.86
$ws_BSM.Range("I134").Value = 891.975
$ws_BSM.Range("J134").Value = 1831.4
$ws_BSM.Range("K134").Value = 2675.925
$ws_BSM.Range("L134").Value = 5494.200000000001
$ws_BSM.Range("M134").Value = -140.9250000000002
$ws_BSM.Range("N134").Value = -10564.2

# CRP row 132
$ws_CRP.Range("H132").Value = 20837644
$ws_CRP.Range("I132").Value = 32263354
$ws_CRP.Range("J132").Value = 2524.353
$ws_CRP.Range("K132").Value = 96790062
$ws_CRP.Range("L132").Value = 7573.059
$ws_CRP.Range("M132").Value = -96787532
$ws_CRP.Range("N132").Value = -12633.059

# CRP row 134
$ws_CRP.Range("H134").Value = 2645.4768
$ws_CRP.Range("I134").Value = 3072.68
$ws_CRP.Range("J134").Value = 1221.4667
$ws_CRP.Range("K134").Value = 9218.039999999999
$ws_CRP.Range("L134").Value = 3664.4001
$ws_CRP.Range("M134").Value = -6683.039999999999
$ws_CRP.Range("N134").Value = -8734.400099999999

# CUL row 133
$ws_CUL.Range("H133").Value = 559569.1
$ws_CUL.Range("I133").Value = 4015.875
$ws_CUL.Range("J133").Value = 1004011.7
$ws_CUL.Range("K133").Value = 12047.625
$ws_CUL.Range("L133").Value = 3012035.1
$ws_CUL.Range("M133").Value = -6987.625
$ws_CUL.Range("N133").Value = -3022155.1

# GSM row 132
$ws_GSM.Range("H132").Value = 3723.2712
$ws_GSM.Range("I132").Value = 4250.524
$ws_GSM.Range("J132").Value = 2420.647
$ws_GSM.Range("K132").Value = 12751.572
$ws_GSM.Range("L132").Value = 7261.941
$ws_GSM.Range("M132").Value = -10221.572
$ws_GSM.Range("N132").Value = -12321.941

# LTW row 93
$ws_LTW.Range("H93").Value = 952.11536
$ws_LTW.Range("I93").Value = 943.5294
$ws_LTW.Range("J93").Value = 968.3333
$ws_LTW.Range("K93").Value = 943.5294
$ws_LTW.Range("L93").Value = 968.3333
$ws_LTW.Range("M93").Value = 304.4706
$ws_LTW.Range("N93").Value = -3464.3333

# WVR row 122
$ws_WVR.Range("H122").Value = 20170202
$ws_WVR.Range("I122").Value = 20001560
$ws_WVR.Range("K122").Value = 60004680
$ws_WVR.Range("M122").Value = -60002230

# WVR row 132
$ws_WVR.Range("H132").Value = 7240.55
$ws_WVR.Range("I132").Value = 11399.728
$ws_WVR.Range("J132").Value = 2157.111
$ws_WVR.Range("K132").Value = 34199.18399999999
$ws_WVR.Range("L132").Value = 6471.333
$ws_WVR.Range("M132").Value = -31669.18399999999
$ws_WVR.Range("N132").Value = -11531.333

# WVR row 136
$ws_WVR.Range("H136").Value = 2168
$ws_WVR.Range("I136").Value = 2831.262
$ws_WVR.Range("J136").Value = 1453.7179
$ws_WVR.Range("K136").Value = 8493.786
$ws_WVR.Range("L136").Value = 4361.153700000001
$ws_WVR.Range("M136").Value = -5943.786
$ws_WVR.Range("N136").Value = -9461.153700000001

# WVR row 141
$ws_WVR.Range("H141").Value = 40157.5
$ws_WVR.Range("J141").Value = 40157.5
$ws_WVR.Range("L141").Value = 40157.5
$ws_WVR.Range("N141").Value = -50517.5

